$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (current UNITE column) to host CATEGORY,
# shifting the existing UNITE column to E.
$ws.Columns.Item(4).Insert()

# Header row (row 3)
$ws.Range("D3").Value = "CATEGORY"
$ws.Range("E3").Value = "UNITE"

# Data rows: set CATEGORY = SUPERMARKET, keep UNITE = kg in column E
$ws.Range("D4").Value = "SUPERMARKET"
$ws.Range("E4").Value = "kg"

$ws.Range("D5").Value = "SUPERMARKET"
$ws.Range("E5").Value = "kg"

$ws.Range("D6").Value = "SUPERMARKET"
$ws.Range("E6").Value = "kg"

# Apply the same style as the other body cells to the new column D & E cells
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"

# Column width for column D (13.56)
$ws.Columns.Item(4).ColumnWidth = 13.56

# Row heights
$ws.Rows.Item(3).RowHeight = 16.15
$ws.Rows.Item(4).RowHeight = 13.8
$ws.Rows.Item(6).RowHeight = 13.8

# Update selection to D7
$ws.Range("D7").Select()
